# Insert 3 new rows of weekly price data at row 593 (top of the "Tuna" block),
# shifting all the existing rows 593:689 down to 596:692.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("593:595").Insert()

# Shared / unchanged metadata columns for the three new rows (same market,
# region, product classification as the rest of this block).
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107011
$categoria   = "Tuna"
$variedad    = "Sin especificar"
$unidad      = "`$/caja 18 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 18

$fecha = 45154

# Row 593 - Especial
$r = 593
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 150
$ws.Cells.Item($r, 14).Value = 28000
$ws.Cells.Item($r, 15).Value = 30000
$ws.Cells.Item($r, 16).Value = 29000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1611
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 594 - Primera
$r = 594
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 90
$ws.Cells.Item($r, 14).Value = 24000
$ws.Cells.Item($r, 15).Value = 25000
$ws.Cells.Item($r, 16).Value = 24500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1361
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 595 - Segunda
$r = 595
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 70
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 20000
$ws.Cells.Item($r, 16).Value = 19000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1056
$ws.Cells.Item($r, 20).Value = $kgUnidad
